$d = $word.ActiveDocument

$pairs = @(
    @("818÷2=", "740÷7="),
    @("288÷6=", "645÷5="),
    @("828÷5=", "111÷4="),
    @("541÷8=", "979÷9="),
    @("230÷4=", "202÷2="),
    @("557÷6=", "502÷8="),
    @("884÷7=", "645÷3="),
    @("375÷3=", "670÷5="),
    @("891÷7=", "557÷8="),
    @("137÷7=", "882÷3="),
    @("793÷9=", "442÷6="),
    @("714÷5=", "501÷4="),
    @("674÷9=", "849÷4="),
    @("161÷2=", "738÷3="),
    @("704÷2=", "292÷7="),
    @("682÷9=", "167÷8="),
    @("185÷4=", "952÷5="),
    @("797÷4=", "826÷4="),
    @("318÷8=", "683÷3="),
    @("706÷9=", "406÷9="),
    @("702÷5=", "132÷7="),
    @("351÷6=", "187÷2="),
    @("784÷4=", "951÷2="),
    @("469÷9=", "753÷3="),
    @("562÷2=", "763÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
